$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold/centered/bordered) from G1 into H1, then set values
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = 'd_prompts'

$ws.Range("H2").Value = 15.1
$ws.Range("H3").Value = 11.93
$ws.Range("H4").Value = 3.55
$ws.Range("H5").Value = 18.8
$ws.Range("H6").Value = 5.83
$ws.Range("H7").Value = 0.03
$ws.Range("H8").Value = 2.1
$ws.Range("H9").Value = '[ 959  569 1357  483  177  200    4  175   27   20   16    0   18    0
    3    1    0    0    0    0    0    0    0    0    0    0    0    0
    0    0    0    0    0    0    0    0]'
$ws.Range("H10").Value = '[831   1 576   0 696 313   1 620  19 576   1 390]'
$ws.Range("H11").Value = '{''C'': 831, ''C#'': 1, ''D'': 576, ''D#'': 0, ''E'': 696, ''F'': 313, ''F#'': 1, ''G'': 620, ''G#'': 19, ''A'': 576, ''A#'': 1, ''B'': 390}'
$ws.Range("H12").Value = 1.47
$ws.Range("H13").Value = '[[211   0 203   0  59  19   0  48   1 104   0 183]
 [  0   0   0   0   0   0   0   0   1   0   0   0]
 [200   1 106   0 180  10   0  36   0  18   0  23]
 [  0   0   0   0   0   0   0   0   0   0   0   0]
 [ 71   0 200   0 144  81   0 110   2  71   0  13]
 [ 12   0   7   0 109 116   0  58   1   7   0   2]
 [  0   0   0   0   1   0   0   0   0   0   0   0]
 [ 50   0  32   0 146  77   0 179   1 109   0  24]
 [  2   0   0   0   0   0   0   1  11   3   1   1]
 [ 91   0  18   0  51   8   1 153   2 158   0  92]
 [  1   0   0   0   0   0   0   0   0   0   0   0]
 [188   0  10   0   6   2   0  30   0 101   0  52]]'
$ws.Range("H14").Value = 268.27
$ws.Range("H15").Value = 6.29
$ws.Range("H16").Value = 2.57
$ws.Range("H17").Value = '[  79 1278 1302  621  380   68   56   32   85   16   14   10   20    8
    9    6   22    7    3    0    2    0    1    0    0    0    0    0
    1    0    0    0]'
$ws.Range("H18").Value = 6.25
$ws.Range("H19").Value = '[   1 1186    0  918    0  243    0 1283    0   56    0   33   12   16
    0  167    0    5   13    7    4    7    0   17    0    4    3    6
    1    3    0   37    0    0    0    0    0    0    0    0    0    0
    0    0    0    0    0    0    0    0    0    0    0    0    0    0
    0    0    0    0    0    0    0    0]'
$ws.Range("H20").Value = '{0: 0, 1: 1, 2: 1186, 3: 0, 4: 918, 5: 0, 6: 243, 7: 0, 8: 1283, 9: 0, 10: 56, 11: 0, 12: 33, 13: 12, 14: 16, 15: 0, 16: 167, 17: 0, 18: 5, 19: 13, 20: 7, 21: 4, 22: 7, 23: 0, 24: 17, 25: 0, 26: 4, 27: 3, 28: 6, 29: 1, 30: 3, 31: 0, 32: 37}'
$ws.Range("H21").Value = '[[  0   1   0 ...   0   0   0]
 [  0 441   0 ...   2   0   4]
 [  0   0   0 ...   0   0   0]
 ...
 [  0   1   0 ...   0   0   0]
 [  0   0   0 ...   0   0   0]
 [  0   9   0 ...   0   0   1]]'
